$wb = $excel.ActiveWorkbook

# Remove the "selected" state of the Classes tab; the new Errors tab will become active
$classesSheet = $wb.Worksheets.Item("Classes")
$classesSheet.Select()

# Add "Errors" sheet after "Classes"
$errorsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $classesSheet)
$errorsSheet.Name = "Errors"
$errorsSheet.Range("A1").Value = 'Sheet "Classes" Row: 2 Missing "OFF CLS"'

# Add "Warnings" sheet after "Errors"
$warningsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $errorsSheet)
$warningsSheet.Name = "Warnings"

# Make Errors sheet the active tab
$errorsSheet.Select()
